$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C ("Förändrad" date) from 45208 to 45212 for all data rows (2-56)
for ($row = 2; $row -le 56; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value = 45212
    }
}
